$wb = $excel.ActiveWorkbook

# Sheet ALC, row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1145.3214
$ws.Range("I15").Value = 1145.3214
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 3435.9642
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -3266.9642

# Sheet ALC, row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 7297.8
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 7297.8
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 7297.8
$ws.Range("N113").Value = -13805.8

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 25003028
$ws.Range("I132").Value = 27030176
$ws.Range("J132").Value = 1535
$ws.Range("K132").Value = 81090528
$ws.Range("L132").Value = 4605
$ws.Range("M132").Value = -81087998
$ws.Range("N132").Value = -9665

# Sheet ALC, row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3999
$ws.Range("I141").Value = 3999
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 11997
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -6817

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3167.8806
$ws.Range("I32").Value = 2376.0508
$ws.Range("J32").Value = 9007.625
$ws.Range("K32").Value = 2376.0508
$ws.Range("L32").Value = 9007.625
$ws.Range("M32").Value = -2089.0508

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2288.889
$ws.Range("I61").Value = 2122.652
$ws.Range("J61").Value = 3244.75
$ws.Range("K61").Value = 2122.652
$ws.Range("L61").Value = 3244.75
$ws.Range("M61").Value = -1910.652

# Sheet ARM, row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2141326.2
$ws.Range("I110").Value = 2530294.8
$ws.Range("J110").Value = 1999.5
$ws.Range("K110").Value = 2530294.8
$ws.Range("L110").Value = 1999.5
$ws.Range("M110").Value = -2528249.8

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 653039
$ws.Range("I122").Value = 1745.1154
$ws.Range("J122").Value = 3475312.5
$ws.Range("K122").Value = 5235.3462
$ws.Range("L122").Value = 10425937.5
$ws.Range("M122").Value = -2785.3462
$ws.Range("N122").Value = -10430837.5

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2288.889
$ws.Range("I136").Value = 2122.652
$ws.Range("J136").Value = 3244.75
$ws.Range("K136").Value = 6367.956
$ws.Range("L136").Value = 9734.25
$ws.Range("M136").Value = -3817.956

# Sheet BSM, row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1242.2593
$ws.Range("I20").Value = 950.8
$ws.Range("J20").Value = 1606.5834
$ws.Range("K20").Value = 950.8
$ws.Range("L20").Value = 1606.5834
$ws.Range("M20").Value = -703.8
$ws.Range("N20").Value = -2100.5834

# Sheet BSM, row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 23813188
$ws.Range("I107").Value = 35715784
$ws.Range("J107").Value = 7999
$ws.Range("K107").Value = 35715784
$ws.Range("L107").Value = 7999
$ws.Range("M107").Value = -35713864
$ws.Range("N107").Value = -11839

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4826.92
$ws.Range("I134").Value = 2266.5
$ws.Range("J134").Value = 8085.636
$ws.Range("K134").Value = 6799.5
$ws.Range("L134").Value = 24256.908
$ws.Range("M134").Value = -4264.5

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28337.383
$ws.Range("I31").Value = 1029.3043
$ws.Range("J31").Value = 85436.09
$ws.Range("K31").Value = 1029.3043
$ws.Range("L31").Value = 85436.09
$ws.Range("M31").Value = -734.3043
$ws.Range("N31").Value = -86026.09

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 28337.383
$ws.Range("I34").Value = 1029.3043
$ws.Range("J34").Value = 85436.09
$ws.Range("K34").Value = 1029.3043
$ws.Range("L34").Value = 85436.09
$ws.Range("M34").Value = -827.3043
$ws.Range("N34").Value = -85840.09

# Sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2016.7407
$ws.Range("I58").Value = 2006.6666
$ws.Range("J58").Value = 2024.8
$ws.Range("K58").Value = 2006.6666
$ws.Range("L58").Value = 2024.8
$ws.Range("M58").Value = -1803.6666
$ws.Range("N58").Value = -2430.8

# Sheet CRP, row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4993.125
$ws.Range("I99").Value = 3737.5
$ws.Range("J99").Value = 6248.75
$ws.Range("K99").Value = 3737.5
$ws.Range("L99").Value = 6248.75
$ws.Range("M99").Value = -2239.5
$ws.Range("N99").Value = -9244.75

# Sheet CRP, row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4993.125
$ws.Range("I126").Value = 3737.5
$ws.Range("J126").Value = 6248.75
$ws.Range("K126").Value = 11212.5
$ws.Range("L126").Value = 18746.25
$ws.Range("M126").Value = -8742.5
$ws.Range("N126").Value = -23686.25

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 65893.92999999999
$ws.Range("I132").Value = 1846.4445
$ws.Range("J132").Value = 181179.4
$ws.Range("K132").Value = 5539.333500000001
$ws.Range("L132").Value = 543538.2
$ws.Range("M132").Value = -3009.333500000001

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 36406.48
$ws.Range("I134").Value = 58705.375
$ws.Range("J134").Value = 3971.7273
$ws.Range("K134").Value = 176116.125
$ws.Range("L134").Value = 11915.1819
$ws.Range("M134").Value = -173581.125

# Sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2016.7407
$ws.Range("I136").Value = 2006.6666
$ws.Range("J136").Value = 2024.8
$ws.Range("K136").Value = 6019.9998
$ws.Range("L136").Value = 6074.4
$ws.Range("M136").Value = -3469.9998
$ws.Range("N136").Value = -11174.4

# Sheet GSM, row 12
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = ""

# Sheet GSM, row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13340217
$ws.Range("I70").Value = 14292661
$ws.Range("J70").Value = 6000
$ws.Range("K70").Value = 14292661
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -14292391
$ws.Range("N70").Value = -6540

# Sheet GSM, row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 13340217
$ws.Range("I73").Value = 14292661
$ws.Range("J73").Value = 6000
$ws.Range("K73").Value = 14292661
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -14291725
$ws.Range("N73").Value = -7872

# Sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 7554700
$ws.Range("I102").Value = 11114010
$ws.Range("J102").Value = 3105563.5
$ws.Range("K102").Value = 11114010
$ws.Range("L102").Value = 3105563.5
$ws.Range("M102").Value = -11112388

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 345237.3
$ws.Range("I122").Value = 425484.28
$ws.Range("J122").Value = 8200
$ws.Range("K122").Value = 1276452.84
$ws.Range("L122").Value = 24600
$ws.Range("M122").Value = -1274002.84

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4139.875
$ws.Range("I132").Value = 3181.7273
$ws.Range("J132").Value = 6247.8
$ws.Range("K132").Value = 9545.1819
$ws.Range("L132").Value = 18743.4
$ws.Range("M132").Value = -7015.1819

# Sheet GSM, row 137
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""

# Sheet LTW, row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3474477.8
$ws.Range("I61").Value = 4117548
$ws.Range("J61").Value = 1898.4
$ws.Range("K61").Value = 4117548
$ws.Range("L61").Value = 1898.4
$ws.Range("M61").Value = -4117346
$ws.Range("N61").Value = -2302.4

# Sheet LTW, row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3474477.8
$ws.Range("I113").Value = 4117548
$ws.Range("J113").Value = 1898.4
$ws.Range("K113").Value = 4117548
$ws.Range("L113").Value = 1898.4
$ws.Range("M113").Value = -4115378
$ws.Range("N113").Value = -6238.4

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5217.9165
$ws.Range("I122").Value = 3546.5715
$ws.Range("J122").Value = 7557.8
$ws.Range("K122").Value = 10639.7145
$ws.Range("L122").Value = 22673.4
$ws.Range("M122").Value = -8189.7145

# Sheet LTW, row 130
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H130").Value = 63500
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 63500
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 63500
$ws.Range("N130").Value = -73540

# Sheet WVR, row 42
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 33999
$ws.Range("I42").Value = 33999
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 33999
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -33621
$ws.Range("N42").Value = ""

# Sheet WVR, row 43
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = ""

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2129.5454
$ws.Range("I122").Value = 1241.25
$ws.Range("J122").Value = 4498.3335
$ws.Range("K122").Value = 3723.75
$ws.Range("L122").Value = 13495.0005
$ws.Range("M122").Value = -1273.75

# Sheet WVR, row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1845.8667
$ws.Range("I126").Value = 1891.6428
$ws.Range("J126").Value = 1205
$ws.Range("K126").Value = 5674.928400000001
$ws.Range("L126").Value = 3615
$ws.Range("M126").Value = -3204.928400000001
$ws.Range("N126").Value = -8555
